# Fixed retail indicators in union:
# Insert two new rows (Id/Name 19, 20) into the "r AnalysisUnit_Variable"
# sheet just above the existing "...IND_26" (index 47) row, shifting the
# remaining rows down by two, and move the active selection/view to the
# newly inserted rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert two blank rows above row 20 (old rows 20-98 shift to 22-100)
$ws.Range("A20:A21").EntireRow.Insert()

# Fill column A (Action = CREATE/MODIFY) for both new rows first
$ws.Cells.Item(20, 1).Value2 = "CREATE/MODIFY"
$ws.Cells.Item(21, 1).Value2 = "CREATE/MODIFY"

# Fill column B (Id) for both new rows - builds shared strings 224/225
$ws.Cells.Item(20, 2).Value2 = "COUNTERPARTY_RETAIL_IND_19"
$ws.Cells.Item(21, 2).Value2 = "COUNTERPARTY_RETAIL_IND_20"

# Fill column C (Name) for both new rows - reuses shared strings 224/225
$ws.Cells.Item(20, 3).Value2 = "COUNTERPARTY_RETAIL_IND_19"
$ws.Cells.Item(21, 3).Value2 = "COUNTERPARTY_RETAIL_IND_20"

# Fill column E (Relation table) for both new rows - reuses existing string
$ws.Cells.Item(20, 5).Value2 = "COUNTERPARTY_RETAIL"
$ws.Cells.Item(21, 5).Value2 = "COUNTERPARTY_RETAIL"

# Fill column F (Formula) for both new rows - builds shared strings 226/227
$ws.Cells.Item(20, 6).Value2 = "RETAIL_IND_19"
$ws.Cells.Item(21, 6).Value2 = "RETAIL_IND_20"

# Move the view/selection onto the newly-inserted rows
$ws.Activate()
$ws.Range("A20:XFD21").Select()
